$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row text to lowercase for consistency
$ws.Range("A1").Value = "question"
$ws.Range("B1").Value = "context"

# Remove the alignment style previously applied to the data rows
# by resetting the style of the whole used range back to Normal.
$usedRange = $ws.UsedRange
$usedRange.Style = "Normal"

# Restore default selection to B1
$ws.Range("B1").Select()
